$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8577154308617234
$ws.Range("B4").Value = 0.8644501278772379
$ws.Range("B5").Value = 0.6293706293706294
$ws.Range("B6").Value = 0.949438202247191
$ws.Range("B7").Value = 0.9186381074168799
